# Applies the cfb_weather.xlsx update described by the commit:
#   "Update cfb_weather.xlsx with Timestamp 2025-10-03T10:01:52.534649"
#
# The workbook has two sheets:
#   FBS   - columns A..AK, AK = Timestamp (all rows share the same stamp)
#   Other - columns A..X (no Timestamp column)

$wb = $excel.ActiveWorkbook

$fbs   = $wb.Worksheets.Item("FBS")
$other = $wb.Worksheets.Item("Other")

# ---------------------------------------------------------------------
# 1. Refresh run Timestamp (every AK cell on FBS carries the same value)
# ---------------------------------------------------------------------
$newTimestamp = "2025-10-03T10:01:52.534649"
$lastRow = $fbs.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $fbs.Cells.Item($r, 37)  # column AK = 37
    if ($cell.Value -ne $null -and $cell.Value -ne "") {
        $cell.Value = $newTimestamp
    }
}

# ---------------------------------------------------------------------
# 2. FBS sheet - updated weather / odds figures
# ---------------------------------------------------------------------
$fbs.Range("Q10").Value = "NW"

$fbs.Range("Y17").Value = 54.5
$fbs.Range("Z17").Value = -105
$fbs.Range("AE17").Value = -0.01801801801801802

$fbs.Range("Q19").Value = "ESE"

$fbs.Range("Q20").Value = "NW"

$fbs.Range("N24").Value = "SW"
$fbs.Range("O24").Value = 78.56
$fbs.Range("P24").Value = 9.6
$fbs.Range("Q24").Value = "SW"
$fbs.Range("U24").Value = 0.5

$fbs.Range("O25").Value = 61.04
$fbs.Range("P25").Value = 5.1
$fbs.Range("Q25").Value = "SSW"
$fbs.Range("U25").Value = -0.9

$fbs.Range("M26").Value = "N"
$fbs.Range("N26").Value = "NW"
$fbs.Range("O26").Value = 69.55999999999999
$fbs.Range("P26").Value = 5.3
$fbs.Range("Q26").Value = "NW"
$fbs.Range("U26").Value = -0.2

$fbs.Range("N27").Value = "NNE"
$fbs.Range("O27").Value = 63.77
$fbs.Range("P27").Value = 3.3
$fbs.Range("U27").Value = -2.9

$fbs.Range("Z28").Value = -115

$fbs.Range("Z30").Value = -106

$fbs.Range("Z34").Value = -115

$fbs.Range("Q36").Value = "N"

$fbs.Range("Z37").Value = -115

$fbs.Range("Q38").Value = "NE"
$fbs.Range("Y38").Value = 48.5
$fbs.Range("Z38").Value = -106
$fbs.Range("AE38").Value = -0.0396039603960396

$fbs.Range("Z46").Value = -118

# ---------------------------------------------------------------------
# 3. Other sheet - updated weather figures
# ---------------------------------------------------------------------
$other.Range("Q16").Value = 61.22
$other.Range("R16").Value = 4.4

$other.Range("Q17").Value = 58.34
$other.Range("R17").Value = 5.6

$other.Range("O18").Value = "NNE"
$other.Range("P18").Value = "N"
$other.Range("Q18").Value = 59.57
$other.Range("R18").Value = 4.7
$other.Range("S18").Value = "N"

$other.Range("S26").Value = "E"

$other.Range("S38").Value = "NW"
